# Partial refactoring of MilestoneData: rename "Assurance MM6*" labels to
# "Approval MM6*" and add "Yes" markers for the IPDC milestones re-baseline row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Assurance MM6" family of labels in column A (rows 2-8) to "Approval MM6".
for ($r = 2; $r -le 8; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $text = $cell.Text
    $cell.Value = $text -replace '^Assurance MM6', 'Approval MM6'
}

# Row 14 ("Re-baseline IPDC milestones") gets "Yes" markers in B:D, matching
# the pattern already used by row 10 ("Re-baseline this quarter" -> "No").
$ws.Range("B14:D14").Value = "Yes"

# Selection moved.
$ws.Range("G25").Select()
